# Welcome.docx edits
#
# 1) Merge " would have played, which " + "may" + " help you to progress."
#    into a single run (text unchanged, only run-splitting is collapsed).
# 2) "...option to enter a<NBSP>nickname to be stored..." ->
#    "...option to enter your<NBSP>first name to be stored..."
#    (the word "nickname" is split as "nick"+"name"; "a<NBSP>nick" becomes
#    "your<NBSP>first ", and the trailing "name..." text is left untouched).
# 3) "(rules, interface usage, ... ) and<NBSP>for (basic) " ->
#    "(interface usage, rules, ... ) and<NBSP>for basic "
#    (swap "rules," after "interface usage," and drop the parens around
#    "basic").
# 4) Merge "C" + "lick on the       button for " into a single run.
# 5) Normal style: suppress automatic hyphenation (adds
#    <w:suppressAutoHyphens/> right after <w:widowControl/>).
#
# NOTE: the source document makes heavy use of non-breaking spaces (U+00A0)
# as word joiners; $nbsp below is used everywhere the original text used one,
# so unaffected text keeps its original joiners.

$d = $word.ActiveDocument
$nbsp = [char]0x00A0

# --- 1) "would have played, which" + "may" + "help you to progress" ---
$d.Content.Find.Execute(
    " would have played, which may help you to progress.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    " would have played, which may help you to progress.",
    2) | Out-Null

# --- 2) "a nickname" -> "your first name" ---
$d.Content.Find.Execute(
    "a${nbsp}nick",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "your${nbsp}first ",
    2) | Out-Null

# --- 3) reorder "rules," / "interface usage," and drop parens around "basic" ---
$d.Content.Find.Execute(
    "(rules, interface usage, ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "(interface usage, rules, ",
    2) | Out-Null

$d.Content.Find.Execute(
    "and${nbsp}for (basic) ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "and${nbsp}for basic ",
    2) | Out-Null

# --- 4) "C" + "lick on the       button for " -> "Click on the       button for " ---
$d.Content.Find.Execute(
    "Click on the       button for ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Click on the       button for ",
    2) | Out-Null

# --- 5) Normal style: suppress automatic hyphenation ---
$normal = $d.Styles("Normal")
$normal.ParagraphFormat.Hyphenation = $false

Write-Host "edits applied"
